$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marks")

# Session 3 (Divide & Conquer) marks entry
$ws.Range("C4").Value2 = 10
$ws.Range("D4").Value2 = 10

# Update comments: C5 (Session 2 comment) simplified, D5 (Session 3 comment) new
$ws.Range("C5").Value2 = "Good work"
$ws.Range("D5").Value2 = "Very good"

# Update the active selection shown in the sheet view
$ws.Range("E5:E12").Select()
